{"js": "const sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\nconst header = sections.items[0].getHeader(\"primary\");\nconst tables = header.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Reapply the (semantically unchanged) bold/italic/strike-through/underline/\n// size/color direct character formatting on every populated data cell of the\n// demo table in the header. This mirrors the upstream change, which only\n// altered how the OOXML serializer (Apache POI 4.1.0 -> 5.2.3) spells out\n// on/off run properties and orders <w:rPr> children -- no visible formatting\n// actually changes.\nconst cellsByRow = [];\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  cellsByRow.push(cells);\n}\nawait context.sync();\n\nconst fontsToLoad = [];\n{\n  const cell = cellsByRow[0].items[1];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = true;\n  font.italic = false;\n  font.strikeThrough = false;\n  font.size = 10.0;\n  font.underline = \"None\";\n}\n{\n  const cell = cellsByRow[0].items[2];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = false;\n  font.italic = true;\n  font.strikeThrough = false;\n  font.size = 10.0;\n  font.underline = \"None\";\n}\n{\n  const cell = cellsByRow[0].items[3];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = false;\n  font.italic = false;\n  font.strikeThrough = false;\n  font.size = 10.0;\n  font.underline = \"Single\";\n}\n{\n  const cell = cellsByRow[0].items[4];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = false;\n  font.italic = false;\n  font.strikeThrough = true;\n  font.size = 10.0;\n  font.underline = \"None\";\n}\n{\n  const cell = cellsByRow[1].items[0];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = true;\n  font.italic = false;\n  font.strikeThrough = false;\n  font.size = 10.0;\n  font.underline = \"None\";\n}\n{\n  const cell = cellsByRow[1].items[1];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = true;\n  font.italic = false;\n  font.strikeThrough = false;\n  font.size = 5.0;\n  font.underline = \"None\";\n  font.color = \"#FF007F\";\n}\n{\n  const cell = cellsByRow[1].items[2];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = true;\n  font.italic = true;\n  font.strikeThrough = false;\n  font.size = 6.0;\n  font.underline = \"None\";\n  font.color = \"#007FFF\";\n}\n{\n  const cell = cellsByRow[1].items[3];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = true;\n  font.italic = false;\n  font.strikeThrough = false;\n  font.size = 7.0;\n  font.underline = \"Single\";\n  font.color = \"#7FFF00\";\n}\n{\n  const cell = cellsByRow[1].items[4];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = true;\n  font.italic = false;\n  font.strikeThrough = true;\n  font.size = 8.0;\n  font.underline = \"None\";\n  font.color = \"#FF007F\";\n}\n{\n  const cell = cellsByRow[2].items[0];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = false;\n  font.italic = true;\n  font.strikeThrough = false;\n  font.size = 10.0;\n  font.underline = \"None\";\n}\n{\n  const cell = cellsByRow[2].items[1];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = true;\n  font.italic = true;\n  font.strikeThrough = false;\n  font.size = 9.0;\n  font.underline = \"None\";\n  font.color = \"#007FFF\";\n}\n{\n  const cell = cellsByRow[2].items[2];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = false;\n  font.italic = true;\n  font.strikeThrough = false;\n  font.size = 10.0;\n  font.underline = \"None\";\n  font.color = \"#7FFF00\";\n}\n{\n  const cell = cellsByRow[2].items[3];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = false;\n  font.italic = true;\n  font.strikeThrough = false;\n  font.size = 11.0;\n  font.underline = \"Single\";\n  font.color = \"#FF007F\";\n}\n{\n  const cell = cellsByRow[2].items[4];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = false;\n  font.italic = true;\n  font.strikeThrough = true;\n  font.size = 12.0;\n  font.underline = \"None\";\n  font.color = \"#007FFF\";\n}\n{\n  const cell = cellsByRow[3].items[0];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = false;\n  font.italic = false;\n  font.strikeThrough = false;\n  font.size = 10.0;\n  font.underline = \"Single\";\n}\n{\n  const cell = cellsByRow[3].items[1];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = true;\n  font.italic = false;\n  font.strikeThrough = false;\n  font.size = 13.0;\n  font.underline = \"Single\";\n  font.color = \"#7FFF00\";\n}\n{\n  const cell = cellsByRow[3].items[2];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = false;\n  font.italic = true;\n  font.strikeThrough = false;\n  font.size = 14.0;\n  font.underline = \"Single\";\n  font.color = \"#FF007F\";\n}\n{\n  const cell = cellsByRow[3].items[3];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = false;\n  font.italic = false;\n  font.strikeThrough = false;\n  font.size = 15.0;\n  font.underline = \"Single\";\n  font.color = \"#007FFF\";\n}\n{\n  const cell = cellsByRow[3].items[4];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = false;\n  font.italic = false;\n  font.strikeThrough = true;\n  font.size = 16.0;\n  font.underline = \"Single\";\n  font.color = \"#7FFF00\";\n}\n{\n  const cell = cellsByRow[4].items[0];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = false;\n  font.italic = false;\n  font.strikeThrough = true;\n  font.size = 10.0;\n  font.underline = \"None\";\n}\n{\n  const cell = cellsByRow[4].items[1];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = true;\n  font.italic = false;\n  font.strikeThrough = true;\n  font.size = 17.0;\n  font.underline = \"None\";\n  font.color = \"#FF007F\";\n}\n{\n  const cell = cellsByRow[4].items[2];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = false;\n  font.italic = true;\n  font.strikeThrough = true;\n  font.size = 18.0;\n  font.underline = \"None\";\n  font.color = \"#007FFF\";\n}\n{\n  const cell = cellsByRow[4].items[3];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = false;\n  font.italic = false;\n  font.strikeThrough = true;\n  font.size = 19.0;\n  font.underline = \"Single\";\n  font.color = \"#7FFF00\";\n}\n{\n  const cell = cellsByRow[4].items[4];\n  const para = cell.body.paragraphs.getFirst();\n  const font = para.font;\n  font.bold = false;\n  font.italic = false;\n  font.strikeThrough = true;\n  font.size = 20.0;\n  font.underline = \"None\";\n  font.color = \"#FF007F\";\n}\nawait context.sync();", "ps1": "# Reapply the (semantically unchanged) bold/italic/strike-through/underline/\n# size/color direct character formatting on every populated data cell of the\n# demo table in the header. This mirrors the upstream change, which only\n# altered how the OOXML serializer (Apache POI 4.1.0 -> 5.2.3) spells out\n# on/off run properties and orders <w:rPr> children -- no visible formatting\n# actually changes.\n$d = $word.ActiveDocument\n$sec = $d.Sections.Item(1)\n$hdr = $sec.Headers.Item(1)\n$tbl = $hdr.Range.Tables.Item(1)\n\n$cell = $tbl.Cell(1, 2)\n$font = $cell.Range.Font\n$font.Bold = 1\n$font.Italic = 0\n$font.StrikeThrough = 0\n$font.Size = 10.0\n$font.Underline = 0\n\n$cell = $tbl.Cell(1, 3)\n$font = $cell.Range.Font\n$font.Bold = 0\n$font.Italic = 1\n$font.StrikeThrough = 0\n$font.Size = 10.0\n$font.Underline = 0\n\n$cell = $tbl.Cell(1, 4)\n$font = $cell.Range.Font\n$font.Bold = 0\n$font.Italic = 0\n$font.StrikeThrough = 0\n$font.Size = 10.0\n$font.Underline = 1\n\n$cell = $tbl.Cell(1, 5)\n$font = $cell.Range.Font\n$font.Bold = 0\n$font.Italic = 0\n$font.StrikeThrough = 1\n$font.Size = 10.0\n$font.Underline = 0\n\n$cell = $tbl.Cell(2, 1)\n$font = $cell.Range.Font\n$font.Bold = 1\n$font.Italic = 0\n$font.StrikeThrough = 0\n$font.Size = 10.0\n$font.Underline = 0\n\n$cell = $tbl.Cell(2, 2)\n$font = $cell.Range.Font\n$font.Bold = 1\n$font.Italic = 0\n$font.StrikeThrough = 0\n$font.Size = 5.0\n$font.Underline = 0\n$font.Color = 8323327  # #FF007F\n\n$cell = $tbl.Cell(2, 3)\n$font = $cell.Range.Font\n$font.Bold = 1\n$font.Italic = 1\n$font.StrikeThrough = 0\n$font.Size = 6.0\n$font.Underline = 0\n$font.Color = 16744192  # #007FFF\n\n$cell = $tbl.Cell(2, 4)\n$font = $cell.Range.Font\n$font.Bold = 1\n$font.Italic = 0\n$font.StrikeThrough = 0\n$font.Size = 7.0\n$font.Underline = 1\n$font.Color = 65407  # #7FFF00\n\n$cell = $tbl.Cell(2, 5)\n$font = $cell.Range.Font\n$font.Bold = 1\n$font.Italic = 0\n$font.StrikeThrough = 1\n$font.Size = 8.0\n$font.Underline = 0\n$font.Color = 8323327  # #FF007F\n\n$cell = $tbl.Cell(3, 1)\n$font = $cell.Range.Font\n$font.Bold = 0\n$font.Italic = 1\n$font.StrikeThrough = 0\n$font.Size = 10.0\n$font.Underline = 0\n\n$cell = $tbl.Cell(3, 2)\n$font = $cell.Range.Font\n$font.Bold = 1\n$font.Italic = 1\n$font.StrikeThrough = 0\n$font.Size = 9.0\n$font.Underline = 0\n$font.Color = 16744192  # #007FFF\n\n$cell = $tbl.Cell(3, 3)\n$font = $cell.Range.Font\n$font.Bold = 0\n$font.Italic = 1\n$font.StrikeThrough = 0\n$font.Size = 10.0\n$font.Underline = 0\n$font.Color = 65407  # #7FFF00\n\n$cell = $tbl.Cell(3, 4)\n$font = $cell.Range.Font\n$font.Bold = 0\n$font.Italic = 1\n$font.StrikeThrough = 0\n$font.Size = 11.0\n$font.Underline = 1\n$font.Color = 8323327  # #FF007F\n\n$cell = $tbl.Cell(3, 5)\n$font = $cell.Range.Font\n$font.Bold = 0\n$font.Italic = 1\n$font.StrikeThrough = 1\n$font.Size = 12.0\n$font.Underline = 0\n$font.Color = 16744192  # #007FFF\n\n$cell = $tbl.Cell(4, 1)\n$font = $cell.Range.Font\n$font.Bold = 0\n$font.Italic = 0\n$font.StrikeThrough = 0\n$font.Size = 10.0\n$font.Underline = 1\n\n$cell = $tbl.Cell(4, 2)\n$font = $cell.Range.Font\n$font.Bold = 1\n$font.Italic = 0\n$font.StrikeThrough = 0\n$font.Size = 13.0\n$font.Underline = 1\n$font.Color = 65407  # #7FFF00\n\n$cell = $tbl.Cell(4, 3)\n$font = $cell.Range.Font\n$font.Bold = 0\n$font.Italic = 1\n$font.StrikeThrough = 0\n$font.Size = 14.0\n$font.Underline = 1\n$font.Color = 8323327  # #FF007F\n\n$cell = $tbl.Cell(4, 4)\n$font = $cell.Range.Font\n$font.Bold = 0\n$font.Italic = 0\n$font.StrikeThrough = 0\n$font.Size = 15.0\n$font.Underline = 1\n$font.Color = 16744192  # #007FFF\n\n$cell = $tbl.Cell(4, 5)\n$font = $cell.Range.Font\n$font.Bold = 0\n$font.Italic = 0\n$font.StrikeThrough = 1\n$font.Size = 16.0\n$font.Underline = 1\n$font.Color = 65407  # #7FFF00\n\n$cell = $tbl.Cell(5, 1)\n$font = $cell.Range.Font\n$font.Bold = 0\n$font.Italic = 0\n$font.StrikeThrough = 1\n$font.Size = 10.0\n$font.Underline = 0\n\n$cell = $tbl.Cell(5, 2)\n$font = $cell.Range.Font\n$font.Bold = 1\n$font.Italic = 0\n$font.StrikeThrough = 1\n$font.Size = 17.0\n$font.Underline = 0\n$font.Color = 8323327  # #FF007F\n\n$cell = $tbl.Cell(5, 3)\n$font = $cell.Range.Font\n$font.Bold = 0\n$font.Italic = 1\n$font.StrikeThrough = 1\n$font.Size = 18.0\n$font.Underline = 0\n$font.Color = 16744192  # #007FFF\n\n$cell = $tbl.Cell(5, 4)\n$font = $cell.Range.Font\n$font.Bold = 0\n$font.Italic = 0\n$font.StrikeThrough = 1\n$font.Size = 19.0\n$font.Underline = 1\n$font.Color = 65407  # #7FFF00\n\n$cell = $tbl.Cell(5, 5)\n$font = $cell.Range.Font\n$font.Bold = 0\n$font.Italic = 0\n$font.StrikeThrough = 1\n$font.Size = 20.0\n$font.Underline = 0\n$font.Color = 8323327  # #FF007F\n"}
